# Commit: "Updated tags - removed MAN4"
#
# The row for BOLOGNA / CASTELFRANCO / MAN4 (tag B025129ABO33_1_J_LIV_FALD004_MFLIV01 /
# B025129ABO33_1_J_POR_MISU004_MFPOR01) is removed entirely from the worksheet,
# shifting every subsequent row up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (1-based) is "BOLOGNA | CASTELFRANCO | MAN4 | ... | 39.5 | 1 | 1".
# Deleting the whole row shifts rows 13-16 up to become rows 12-15.
$ws.Rows.Item(12).Delete()

# Reflect the row selection left by the edit (row 12, now occupied by what used
# to be row 13 - the BoPan8 entry).
$ws.Rows.Item(12).Select()
